# edit.ps1
# Applies the restructuring described in the diff:
#  - Reorders existing monthly rows so that Oct/Nov/Dec of each year are
#    moved to the front of that year's block (immediately preceding Jan)
#  - Appends new data rows for 2022-10 .. 2023-07
#  - Expands the sheet dimension from A1:C49 to A1:C68

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 67,3

$data[0,0] = "2018-10"
$data[0,1] = [double]117.9
$data[0,2] = [double]105.5
$data[1,0] = "2018-11"
$data[1,1] = [double]117.3
$data[1,2] = [double]106.5
$data[2,0] = "2018-12"
$data[2,1] = [double]110.9
$data[2,2] = [double]105.6
$data[3,0] = "2018-01"
$data[3,1] = [double]118.2217
$data[3,2] = [double]105.6267
$data[4,0] = "2018-02"
$data[4,1] = [double]115.4
$data[4,2] = [double]105.3
$data[5,0] = "2018-03"
$data[5,1] = [double]114
$data[5,2] = [double]103.7
$data[6,0] = "2018-04"
$data[6,1] = [double]113.6
$data[6,2] = [double]104.7
$data[7,0] = "2018-05"
$data[7,1] = [double]117.8
$data[7,2] = [double]105.9
$data[8,0] = "2018-06"
$data[8,1] = [double]119.6
$data[8,2] = [double]105.4
$data[9,0] = "2018-07"
$data[9,1] = [double]118.5
$data[9,2] = [double]105.1
$data[10,0] = "2018-08"
$data[10,1] = [double]118.3
$data[10,2] = [double]104
$data[11,0] = "2018-09"
$data[11,1] = [double]118.3
$data[11,2] = [double]106.3
$data[12,0] = "2019-10"
$data[12,1] = [double]101.1
$data[12,2] = [double]101.4
$data[13,0] = "2019-11"
$data[13,1] = [double]100.2
$data[13,2] = [double]102.3
$data[14,0] = "2019-12"
$data[14,1] = [double]102.8
$data[14,2] = [double]103
$data[15,0] = "2019-01"
$data[15,1] = [double]105
$data[15,2] = [double]105.7
$data[16,0] = "2019-02"
$data[16,1] = [double]106.7
$data[16,2] = [double]104.5
$data[17,0] = "2019-03"
$data[17,1] = [double]106.7
$data[17,2] = [double]105.1
$data[18,0] = "2019-04"
$data[18,1] = [double]106.7
$data[18,2] = [double]104
$data[19,0] = "2019-05"
$data[19,1] = [double]106.1
$data[19,2] = [double]103.5
$data[20,0] = "2019-06"
$data[20,1] = [double]105.3
$data[20,2] = [double]102.8
$data[21,0] = "2019-07"
$data[21,1] = [double]105.6
$data[21,2] = [double]100.6
$data[22,0] = "2019-08"
$data[22,1] = [double]103.5
$data[22,2] = [double]100.9
$data[23,0] = "2019-09"
$data[23,1] = [double]101.9
$data[23,2] = [double]102.4
$data[24,0] = "2020-10"
$data[24,1] = [double]102
$data[24,2] = [double]101.9
$data[25,0] = "2020-11"
$data[25,1] = [double]102.1
$data[25,2] = [double]101.6
$data[26,0] = "2020-12"
$data[26,1] = [double]104
$data[26,2] = [double]103.4
$data[27,0] = "2020-01"
$data[27,1] = [double]104.4
$data[27,2] = [double]103.6
$data[28,0] = "2020-02"
$data[28,1] = [double]103.3
$data[28,2] = [double]104.1
$data[29,0] = "2020-03"
$data[29,1] = [double]97.59999999999999
$data[29,2] = [double]105.2
$data[30,0] = "2020-04"
$data[30,1] = [double]94.09999999999999
$data[30,2] = [double]103
$data[31,0] = "2020-05"
$data[31,1] = [double]93.40000000000001
$data[31,2] = [double]100.4
$data[32,0] = "2020-06"
$data[32,1] = [double]95.3
$data[32,2] = [double]99.5
$data[33,0] = "2020-07"
$data[33,1] = [double]99.2
$data[33,2] = [double]100.8
$data[34,0] = "2020-08"
$data[34,1] = [double]100.6
$data[34,2] = [double]102.8
$data[35,0] = "2020-09"
$data[35,1] = [double]102.2
$data[35,2] = [double]100.2
$data[36,0] = "2021-10"
$data[36,1] = [double]124.1
$data[36,2] = [double]105.5
$data[37,0] = "2021-11"
$data[37,1] = [double]122.7
$data[37,2] = [double]107.2
$data[38,0] = "2021-12"
$data[38,1] = [double]122.8
$data[38,2] = [double]106
$data[39,0] = "2021-01"
$data[39,1] = [double]103.6
$data[39,2] = [double]101.7
$data[40,0] = "2021-02"
$data[40,1] = [double]109.5
$data[40,2] = [double]102.2
$data[41,0] = "2021-03"
$data[41,1] = [double]117.9
$data[41,2] = [double]104.4
$data[42,0] = "2021-04"
$data[42,1] = [double]124.1
$data[42,2] = [double]106.8
$data[43,0] = "2021-05"
$data[43,1] = [double]128.5
$data[43,2] = [double]109.2
$data[44,0] = "2021-06"
$data[44,1] = [double]125.7
$data[44,2] = [double]106.4
$data[45,0] = "2021-07"
$data[45,1] = [double]122.2
$data[45,2] = [double]105.6
$data[46,0] = "2021-08"
$data[46,1] = [double]121.8
$data[46,2] = [double]104.3
$data[47,0] = "2021-09"
$data[47,1] = [double]123.6
$data[47,2] = [double]105.5
$data[48,0] = "2022-10"
$data[48,1] = [double]90.5
$data[48,2] = [double]102
$data[49,0] = "2022-11"
$data[49,1] = [double]89.5
$data[49,2] = [double]99
$data[50,0] = "2022-12"
$data[50,1] = [double]89.3
$data[50,2] = [double]97.5
$data[51,0] = "2022-01"
$data[51,1] = [double]120.4
$data[51,2] = [double]106.4
$data[52,0] = "2022-02"
$data[52,1] = [double]119.6
$data[52,2] = [double]106.1
$data[53,0] = "2022-03"
$data[53,1] = [double]116.3
$data[53,2] = [double]105.7
$data[54,0] = "2022-04"
$data[54,1] = [double]116
$data[54,2] = [double]105.7
$data[55,0] = "2022-05"
$data[55,1] = [double]109.2
$data[55,2] = [double]104.4
$data[56,0] = "2022-06"
$data[56,1] = [double]107.5
$data[56,2] = [double]107.7
$data[57,0] = "2022-07"
$data[57,1] = [double]97
$data[57,2] = [double]108.2
$data[58,0] = "2022-08"
$data[58,1] = [double]94.2
$data[58,2] = [double]104.4
$data[59,0] = "2022-09"
$data[59,1] = [double]92
$data[59,2] = [double]103.9
$data[60,0] = "2023-01"
$data[60,1] = [double]88.2
$data[60,2] = [double]95.3
$data[61,0] = "2023-02"
$data[61,1] = [double]89.3
$data[61,2] = [double]95.90000000000001
$data[62,0] = "2023-03"
$data[62,1] = [double]88
$data[62,2] = [double]95
$data[63,0] = "2023-04"
$data[63,1] = [double]85.90000000000001
$data[63,2] = [double]92.40000000000001
$data[64,0] = "2023-05"
$data[64,1] = [double]83.7
$data[64,2] = [double]90.8
$data[65,0] = "2023-06"
$data[65,1] = [double]84.09999999999999
$data[65,2] = [double]89.59999999999999
$data[66,0] = "2023-07"
$data[66,1] = [double]93.7
$data[66,2] = [double]89

# Write all data (A2:C68) in a single bulk assignment
$ws.Range("A2:C68").Value = $data

# The newly added rows (50..68) need the same style as the other
# column-A cells (bold, centered, bordered) - copy format down from A49
$ws.Range("A49").Copy()
$ws.Range("A50:A68").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

